# Apply the "Add files via upload" edit to PGWBS.xlsx
# - adds two new task rows (row 7: "login maken", row 8: "database ")
#   to the "project WBS monitor" sheet, both done by "mohamed",
#   each counting 1 hour in column G (so the SUBTOTAL in G12 goes 3 -> 5)
# - moves the active selection on that sheet from C6 to G9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project WBS monitor")

# Fill in the new rows. Set the C column (task name) cells first so the
# shared-string table picks up "login maken" / "database " before "mohamed",
# matching the order new strings were appended to sharedStrings.xml.
$ws.Range("C7").Value = "login maken"
$ws.Range("C8").Value = "database "

$ws.Range("E7").Value = "mohamed"
$ws.Range("E8").Value = "mohamed"

$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1

# Update the active cell/selection on the sheet to G9
$ws.Activate()
$ws.Range("G9").Select()
